$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "debris"
$ws.Range("B18").Value = "/debri:/"
$ws.Range("E18").Value = "s' is silent."
$ws.Range("D18").Value = "U."

# E18 picks up the "quote prefix" formatting used by E17 (e.g. leading apostrophe style)
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("D19").Select()
